# Adds the new "intervention_type" column (K) with a value per clinical trial row,
# as described in the commit: "Add new indicator about 10 years and Drug and list of
# all indicators of all sponsors".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, thin border, centered alignment) from the existing
# header cell J1 onto the new header cell K1, matching the style used by A1:J1.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New header label
$ws.Range("K1").Value = "intervention_type"

# intervention_type value for each clinical trial row (2-74); rows not present in this
# hashtable had no intervention type reported and are left blank.
$interventionTypes = @{
    2 = "RADIATION"
    3 = "DRUG"
    4 = "DRUG"
    5 = "DRUG"
    6 = "DRUG"
    7 = "DRUG"
    8 = "PROCEDURE"
    9 = "DRUG"
    10 = "DRUG"
    11 = "DEVICE"
    12 = "DRUG"
    13 = "RADIATION"
    14 = "DRUG"
    15 = "OTHER"
    16 = "DEVICE"
    17 = "RADIATION"
    18 = "PROCEDURE"
    19 = "BIOLOGICAL"
    20 = "OTHER"
    21 = "DRUG"
    22 = "DRUG"
    23 = "PROCEDURE"
    24 = "DRUG"
    25 = "DRUG"
    26 = "DRUG"
    27 = "DRUG"
    28 = "DRUG"
    30 = "PROCEDURE"
    31 = "DRUG"
    32 = "DRUG"
    33 = "OTHER"
    34 = "PROCEDURE"
    35 = "OTHER"
    36 = "RADIATION"
    38 = "DEVICE"
    39 = "PROCEDURE"
    40 = "DRUG"
    41 = "DRUG"
    42 = "RADIATION"
    43 = "DRUG"
    44 = "DRUG"
    45 = "DRUG"
    46 = "OTHER"
    47 = "OTHER"
    48 = "PROCEDURE"
    49 = "OTHER"
    50 = "DRUG"
    51 = "DRUG"
    52 = "DRUG"
    53 = "DRUG"
    54 = "DRUG"
    55 = "DRUG"
    56 = "OTHER"
    57 = "DRUG"
    58 = "DRUG"
    59 = "RADIATION"
    61 = "DRUG"
    62 = "DRUG"
    63 = "RADIATION"
    64 = "PROCEDURE"
    65 = "DRUG"
    66 = "PROCEDURE"
    67 = "BIOLOGICAL"
    68 = "PROCEDURE"
    69 = "DRUG"
}

foreach ($row in 2..74) {
    $cell = $ws.Cells.Item($row, 11)
    if ($interventionTypes.ContainsKey($row)) {
        $cell.Value = $interventionTypes[$row]
    } else {
        # Materialize an empty cell (matching the blank inline-string cells used
        # elsewhere in this sheet for missing values) without introducing a new style.
        $cell.Value = ""
        $isBold = $cell.Font.Bold
        $cell.Font.Bold = $isBold
    }
}

Write-Output "Done. Used range: $($ws.UsedRange.Address())"
